{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// precedes them (right after the \"LOQ4003: ...\" requirement line).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4003: ...\" requirement paragraph; the blank paragraph and\n// the two footer paragraphs immediately follow it.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4003\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4003' requirement paragraph.\");\n}\n\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i < items.length; i++) {\n  const t = items[i].text;\n  if (\n    t.trim() === \"\" ||\n    t.indexOf(\"Ver no Jupiter\") !== -1 ||\n    t.indexOf(\"Powered by Jekyll\") !== -1\n  ) {\n    toDelete.push(i);\n    // Stop once we've collected the blank line + the two footer lines.\n    if (toDelete.length === 3) break;\n  } else {\n    break;\n  }\n}\n\n// Delete from the highest index down so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  items[toDelete[i]].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n# precedes them (right after the \"LOQ4003: ...\" requirement line).\n$d = $word.ActiveDocument\n\n# Locate the \"LOQ4003: ...\" requirement paragraph.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOQ4003\")\nif (-not $found) {\n    throw \"Could not find the 'LOQ4003' requirement paragraph.\"\n}\n$anchorIndex = $rng.Paragraphs.Item(1).Index\n\n# The blank paragraph and the two footer paragraphs immediately follow it.\n# Delete from the last one back to the first so earlier indices stay valid.\nfor ($k = 3; $k -ge 1; $k--) {\n    $target = $anchorIndex + $k\n    $p = $d.Paragraphs.Item($target)\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"\" -or $text -like \"Ver no Jupiter*\" -or $text -like \"*Powered by Jekyll*\") {\n        $p.Range.Delete()\n    }\n}\n"}
